$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.436.63"
Set-TextValue $ws.Range("E2") "  -2.30%  "
Set-TextValue $ws.Range("D3") "3.173.57"
Set-TextValue $ws.Range("E3") "  -4.25%  "
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "586.36"
Set-TextValue $ws.Range("E5") "  -2.49%  "
Set-TextValue $ws.Range("D6") "134.38"
Set-TextValue $ws.Range("E6") "  -6.01%  "
Set-TextValue $ws.Range("E7") "  -0.31%  "
Set-TextValue $ws.Range("D8") "3.172.76"
Set-TextValue $ws.Range("E8") "  -4.15%  "
Set-TextValue $ws.Range("D9") "0.502"
Set-TextValue $ws.Range("E9") "  -3.71%  "
Set-TextValue $ws.Range("D10") "0.141"
Set-TextValue $ws.Range("E10") "  -5.83%  "
Set-TextValue $ws.Range("E11") "  -5.98%  "
Set-TextValue $ws.Range("D12") "0.450"
Set-TextValue $ws.Range("E12") "  -4.99%  "
Set-TextValue $ws.Range("E13") "  -6.22%  "
Set-TextValue $ws.Range("D14") "32.99"
Set-TextValue $ws.Range("E14") "  -5.30%  "
Set-TextValue $ws.Range("D15") "3.698.13"
Set-TextValue $ws.Range("E15") "  -4.34%  "
Set-TextValue $ws.Range("E16") "  -1.84%  "
Set-TextValue $ws.Range("D17") "3.171.97"
Set-TextValue $ws.Range("E17") "  -4.37%  "
Set-TextValue $ws.Range("D18") "62.421.66"
Set-TextValue $ws.Range("E18") "  -2.50%  "
Set-TextValue $ws.Range("D19") "6.53"
Set-TextValue $ws.Range("E19") "  -5.23%  "
Set-TextValue $ws.Range("D20") "455.81"
Set-TextValue $ws.Range("E20") "  -5.32%  "
Set-TextValue $ws.Range("D21") "13.81"
Set-TextValue $ws.Range("E21") "  -3.07%  "
Set-TextValue $ws.Range("D22") "0.701"
Set-TextValue $ws.Range("E22") "  -4.84%  "
Set-TextValue $ws.Range("D23") "7.61"
Set-TextValue $ws.Range("E23") "  -4.76%  "
Set-TextValue $ws.Range("D24") "13.34"
Set-TextValue $ws.Range("E24") "  -1.84%  "
Set-TextValue $ws.Range("D25") "82.18"
Set-TextValue $ws.Range("E25") "  -2.89%  "
Set-TextValue $ws.Range("E26") "  +0.08%  "
Set-TextValue $ws.Range("E27") "  -0.11%  "
Set-TextValue $ws.Range("D28") "2.67"
Set-TextValue $ws.Range("E28") "  -3.78%  "
Set-TextValue $ws.Range("D29") "6.90"
Set-TextValue $ws.Range("E29") "  -6.09%  "
Set-TextValue $ws.Range("D30") "7.79"
Set-TextValue $ws.Range("E30") "  -4.69%  "
Set-TextValue $ws.Range("E31") "  -7.49%  "
Set-TextValue $ws.Range("D32") "27.16"
Set-TextValue $ws.Range("E32") "  -7.26%  "
Set-TextValue $ws.Range("D33") "0.102"
Set-TextValue $ws.Range("E33") "  -3.86%  "
Set-TextValue $ws.Range("D34") "2.38"
Set-TextValue $ws.Range("E34") "  -6.88%  "
Set-TextValue $ws.Range("E35") "  -6.75%  "
Set-TextValue $ws.Range("D36") "5.77"
Set-TextValue $ws.Range("E36") "  -3.53%  "
Set-TextValue $ws.Range("D37") "51.07"
Set-TextValue $ws.Range("E37") "  -3.26%  "
Set-TextValue $ws.Range("D38") "0.0₃0689"
Set-TextValue $ws.Range("E38") "  -8.58%  "
Set-TextValue $ws.Range("D39") "0.0383"
Set-TextValue $ws.Range("E39") "  -5.43%  "
Set-TextValue $ws.Range("D40") "410.99"
Set-TextValue $ws.Range("E40") "  -4.83%  "
Set-TextValue $ws.Range("D41") "2.932.58"
Set-TextValue $ws.Range("E41") "  -4.06%  "
Set-TextValue $ws.Range("E42") "  -0.25%  "
Set-TextValue $ws.Range("D43") "7.97"
Set-TextValue $ws.Range("E43") "  -5.50%  "
Set-TextValue $ws.Range("D44") "2.62"
Set-TextValue $ws.Range("E44") "  -5.14%  "
Set-TextValue $ws.Range("D45") "0.249"
Set-TextValue $ws.Range("E45") "  -6.81%  "
Set-TextValue $ws.Range("E46") "  -0.04%  "
Set-TextValue $ws.Range("E47") "  -3.85%  "
Set-TextValue $ws.Range("D48") "35.56"
Set-TextValue $ws.Range("E48") "  -1.54%  "
Set-TextValue $ws.Range("D49") "124.60"
Set-TextValue $ws.Range("E49") "  +0.73%  "
Set-TextValue $ws.Range("D50") "25.27"
Set-TextValue $ws.Range("E50") "  -4.61%  "
Set-TextValue $ws.Range("E51") "  -4.26%  "
